{"js": "// Office.js (Word JavaScript API) edit script.\n// 1) Strip the \"Heading2\" paragraph style from the four section headings\n//    (Introduction, Immediate Effects, Long-Term Consequences, Conclusion)\n//    so they fall back to the document default (\"Normal\") styling.\n// 2) Replace the in-text parenthetical citations with their new\n//    reference-id / page-numbered forms (citation check pass):\n//      (Garc\u00eda et al.)        -> (Ref-u782615)     [x2]\n//      (Bishir et al.)        -> (Ref-f173124)\n//      (Newbury et al.)       -> (Ref-f173124)\n//      (Khan and Al-Jahdali)  -> (Pearse et al. 117)\n//      (Hudson et al.)        -> (Pearse et al. 117)\n\nconst body = context.document.body;\n\n// --- 1) Remove the Heading2 style from the section heading paragraphs ---\nconst headings = [\"Introduction\", \"Immediate Effects\", \"Long-Term Consequences\", \"Conclusion\"];\nfor (const heading of headings) {\n  const found = body.search(heading, { matchCase: true, matchWholeWord: true });\n  found.load(\"items,text,style\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    const r = found.items[i];\n    if (r.style === \"Heading 2\") {\n      r.style = \"Normal\";\n    }\n  }\n  await context.sync();\n}\n\n// --- 2) Swap the old author-citations for the new reference codes ---\nconst citationReplacements = [\n  [\"(Garc\u00eda et al.)\", \"(Ref-u782615)\"],\n  [\"(Bishir et al.)\", \"(Ref-f173124)\"],\n  [\"(Newbury et al.)\", \"(Ref-f173124)\"],\n  [\"(Khan and Al-Jahdali)\", \"(Pearse et al. 117)\"],\n  [\"(Hudson et al.)\", \"(Pearse et al. 117)\"],\n];\n\nfor (const [oldText, newText] of citationReplacements) {\n  const hits = body.search(oldText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < hits.items.length; i++) {\n    hits.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# 1) Strip the \"Heading2\" paragraph style from the four section headings\n#    (Introduction, Immediate Effects, Long-Term Consequences, Conclusion)\n#    so they fall back to the document default (\"Normal\") styling.\n# 2) Replace the in-text parenthetical citations with their new\n#    reference-id / page-numbered forms (citation check pass):\n#      (Garc\u00eda et al.)        -> (Ref-u782615)     [x2]\n#      (Bishir et al.)        -> (Ref-f173124)\n#      (Newbury et al.)       -> (Ref-f173124)\n#      (Khan and Al-Jahdali)  -> (Pearse et al. 117)\n#      (Hudson et al.)        -> (Pearse et al. 117)\n\n$d = $word.ActiveDocument\n\n# --- 1) Remove the Heading2 style from the section heading paragraphs ---\n$headingTexts = @(\"Introduction\", \"Immediate Effects\", \"Long-Term Consequences\", \"Conclusion\")\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($headingTexts -contains $t) {\n        if ($p.Style.NameLocal -eq \"Heading 2\") {\n            $p.Style = \"Normal\"\n        }\n    }\n}\n\n# --- 2) Swap the old author-citations for the new reference codes ---\n$wdReplaceAll = 2\n\n$citationReplacements = @(\n    @(\"(Garc\u00eda et al.)\", \"(Ref-u782615)\"),\n    @(\"(Bishir et al.)\", \"(Ref-f173124)\"),\n    @(\"(Newbury et al.)\", \"(Ref-f173124)\"),\n    @(\"(Khan and Al-Jahdali)\", \"(Pearse et al. 117)\"),\n    @(\"(Hudson et al.)\", \"(Pearse et al. 117)\")\n)\n\nforeach ($pair in $citationReplacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, $wdReplaceAll)\n}\n"}
